# "slide on CSS Selectors" - adds a new "CSS Selectors" slide at the end of
# the deck (same Title+Content layout as the other content slides).

$p = $ppt.ActivePresentation

# Re-use the "Title and Content" layout already used by slide 2 ("Goals").
$layout = $p.Slides.Item(2).CustomLayout

# Insert as slide 10 (the new, final slide in the deck).
$newSlide = $p.Slides.AddSlide(10, $layout)

# --- Title -------------------------------------------------------------
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "CSS Selectors"

# --- Body placeholder: position/size override (matches authored layout) --
$body = $newSlide.Shapes.Item(2)
$body.Left = [double]282741 / 12700.0
$body.Top = [double]1825625 / 12700.0
$body.Width = [double]11694695 / 12700.0
$body.Height = [double]4845886 / 12700.0

$tr = $body.TextFrame.TextRange
$tr.Text = "When we need to interact with a web page with Enzyme, need to specify the HTML elements `reg, which button to click, and which text area to fill`r2 main languages to select elements in HTML: CSS Selectors and XPath`rCSS Selectors are the same as when writing .css files`r“.foo”: select all HTML elements with class attribute “foo”`r“#foo”: select element with id attribute “foo” `retc."

# Paragraph 1 (level 1): "When we need to interact with a web page with
# *Enzyme*, need to specify the HTML elements "
$para = $tr.Paragraphs(1)
$r = $para.Characters(47, 6)
$r.Font.Italic = $true

# Paragraph 2 (level 2): "eg, which button to click, and which text area to
# fill"
$para = $tr.Paragraphs(2)
$para.IndentLevel = 2

# Paragraph 3 (level 1): "2 main languages to select elements in HTML: *CSS
# Selectors* and *XPath*"
$para = $tr.Paragraphs(3)
$r = $para.Characters(46, 14)
$r.Font.Italic = $true
$r = $para.Characters(64, 5)
$r.Font.Italic = $true

# Paragraph 4 (level 1): "*CSS Selectors* are the same as when writing *.*
# *css* files"
$para = $tr.Paragraphs(4)
$r = $para.Characters(1, 13)
$r.Font.Italic = $true
$r = $para.Characters(44, 1)
$r.Font.Italic = $true
$r = $para.Characters(45, 3)
$r.Font.Italic = $true

# Paragraph 5 (level 2): "*“.foo”*: select all HTML elements with *class*
# attribute “foo”"
$para = $tr.Paragraphs(5)
$para.IndentLevel = 2
$r = $para.Characters(1, 6)
$r.Font.Italic = $true
$r = $para.Characters(39, 5)
$r.Font.Italic = $true

# Paragraph 6 (level 2): "*“#foo”*: select element with *id* attribute
# “foo” "
$para = $tr.Paragraphs(6)
$para.IndentLevel = 2
$r = $para.Characters(1, 6)
$r.Font.Italic = $true
$r = $para.Characters(29, 2)
$r.Font.Italic = $true

# Paragraph 7 (level 2): "etc."
$para = $tr.Paragraphs(7)
$para.IndentLevel = 2
